# The presentation currently uses the "Integral" (Red Violet) theme colours
# on its one-and-only slide master (ppt/theme/theme1.xml). The commit swaps
# that palette for the stock "Office Theme" colour scheme (the colours that,
# in this deck, had only been living - unused by any slide - in
# ppt/theme/theme2.xml, which backs the notes master).
#
# Re-apply the built-in "Office" theme colours onto the deck's theme via the
# PowerPoint object model's ThemeColorScheme, which is the supported,
# non-destructive way to restyle a design's 12-slot colour scheme
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) from COM automation.

function Convert-RGBToOle {
    param($R, $G, $B)
    return $R + ($G * 256) + ($B * 65536)
}

# Office theme colour scheme (the built-in default PowerPoint palette),
# in clrScheme document order.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

$p = $ppt.ActivePresentation

# Any slide exposes the presentation-wide theme colour scheme; use the
# first one to restyle the deck's single design/master.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $rgb = $officeColors[$i]
    $ole = Convert-RGBToOle $rgb[0] $rgb[1] $rgb[2]
    $themeColors.Colors($i + 1).RGB = $ole
}
